$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1558
$ws1.Range("F5").Value = 161
$ws1.Range("F8").Value = 176
$ws1.Range("F9").Value = 756
$ws1.Range("F12").Value = 344
$ws1.Range("F14").Value = 4
$ws1.Range("F15").Value = 16
$ws1.Range("F16").Value = 6499
$ws1.Range("F20").Value = 162
$ws1.Range("F22").Value = 15521
$ws1.Range("F23").Value = 1538
$ws1.Range("F24").Value = 292
$ws1.Range("F26").Value = 105
$ws1.Range("F27").Value = 11108
$ws1.Range("F28").Value = 768
$ws1.Range("F29").Value = 4351
$ws1.Range("F30").Value = 254
$ws1.Range("F33").Value = 308
$ws1.Range("F34").Value = 128

# Sheet "全部类型" (All types) - same kind of refresh, different row alignment
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1558
$ws4.Range("F5").Value = 161
$ws4.Range("F9").Value = 176
$ws4.Range("F10").Value = 756
$ws4.Range("F14").Value = 344
$ws4.Range("F16").Value = 4
$ws4.Range("F18").Value = 16
$ws4.Range("F19").Value = 6499
$ws4.Range("F23").Value = 162
$ws4.Range("F26").Value = 15521
$ws4.Range("F27").Value = 1538
$ws4.Range("F28").Value = 292
$ws4.Range("F30").Value = 105
$ws4.Range("F32").Value = 11108
$ws4.Range("F33").Value = 768
$ws4.Range("F34").Value = 4351
$ws4.Range("F35").Value = 254
$ws4.Range("F38").Value = 308
$ws4.Range("F39").Value = 128
